$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 68.083336
$ws.Range("I53").Value = 40.25
$ws.Range("J53").Value = 82
$ws.Range("K53").Value = 40.25
$ws.Range("L53").Value = 82
$ws.Range("M53").Value = 596.75
$ws.Range("N53").Value = -1356

$ws.Range("H70").Value = 48860.617
$ws.Range("I70").Value = 167661.67
$ws.Range("J70").Value = 1340.2
$ws.Range("K70").Value = 502985.01
$ws.Range("L70").Value = 4020.6
$ws.Range("M70").Value = -502715.01
$ws.Range("N70").Value = -4560.6

$ws.Range("H73").Value = 48860.617
$ws.Range("I73").Value = 167661.67
$ws.Range("J73").Value = 1340.2
$ws.Range("K73").Value = 502985.01
$ws.Range("L73").Value = 4020.6
$ws.Range("M73").Value = -502049.01
$ws.Range("N73").Value = -5892.6

$ws.Range("H132").Value = 2018.6558
$ws.Range("I132").Value = 1587.8049
$ws.Range("J132").Value = 2901.9
$ws.Range("K132").Value = 4763.4147
$ws.Range("L132").Value = 8705.700000000001
$ws.Range("M132").Value = -2233.4147
$ws.Range("N132").Value = -13765.7

$ws.Range("H135").Value = 773.7586
$ws.Range("I135").Value = 742.1852
$ws.Range("J135").Value = 1200
$ws.Range("K135").Value = 6679.6668
$ws.Range("L135").Value = 10800
$ws.Range("M135").Value = -4144.6668

$ws.Range("H137").Value = 1321.3529
$ws.Range("I137").Value = 979.3214
$ws.Range("J137").Value = 2917.5
$ws.Range("K137").Value = 2937.9642
$ws.Range("L137").Value = 8752.5
$ws.Range("M137").Value = -387.9642000000003
$ws.Range("N137").Value = -13852.5

$ws.Range("H138").Value = 2262.9788
$ws.Range("I138").Value = 1344.5636
$ws.Range("J138").Value = 3558.1794
$ws.Range("K138").Value = 4033.6908
$ws.Range("L138").Value = 10674.5382
$ws.Range("M138").Value = 1106.3092
$ws.Range("N138").Value = -20954.5382

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4954.067
$ws.Range("I61").Value = 4564.125
$ws.Range("J61").Value = 5399.7144
$ws.Range("K61").Value = 4564.125
$ws.Range("L61").Value = 5399.7144
$ws.Range("M61").Value = -4352.125
$ws.Range("N61").Value = -5823.7144

$ws.Range("H63").Value = 5861
$ws.Range("I63").Value = 4826.25
$ws.Range("J63").Value = 10000
$ws.Range("K63").Value = 4826.25
$ws.Range("L63").Value = 10000
$ws.Range("M63").Value = -4140.25
$ws.Range("N63").Value = -11372

$ws.Range("H66").Value = 5861
$ws.Range("I66").Value = 4826.25
$ws.Range("J66").Value = 10000
$ws.Range("K66").Value = 24131.25
$ws.Range("L66").Value = 50000
$ws.Range("M66").Value = -20699.25
$ws.Range("N66").Value = -56864

$ws.Range("H74").Value = 1103.9714
$ws.Range("I74").Value = 1431.7894
$ws.Range("J74").Value = 714.6875
$ws.Range("K74").Value = 1431.7894
$ws.Range("L74").Value = 714.6875
$ws.Range("M74").Value = -557.7893999999999
$ws.Range("N74").Value = -2462.6875

$ws.Range("H77").Value = 1103.9714
$ws.Range("I77").Value = 1431.7894
$ws.Range("J77").Value = 714.6875
$ws.Range("K77").Value = 7158.946999999999
$ws.Range("L77").Value = 3573.4375
$ws.Range("M77").Value = -2790.946999999999
$ws.Range("N77").Value = -12309.4375

$ws.Range("H136").Value = 4954.067
$ws.Range("I136").Value = 4564.125
$ws.Range("J136").Value = 5399.7144
$ws.Range("K136").Value = 13692.375
$ws.Range("L136").Value = 16199.1432
$ws.Range("M136").Value = -11142.375
$ws.Range("N136").Value = -21299.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 331043.25
$ws.Range("I132").Value = 436959.3
$ws.Range("J132").Value = 2703.4
$ws.Range("K132").Value = 1310877.9
$ws.Range("L132").Value = 8110.200000000001
$ws.Range("M132").Value = -1308347.9
$ws.Range("N132").Value = -13170.2

$ws.Range("H141").Value = 48749.75
$ws.Range("I141").Value = 20000
$ws.Range("J141").Value = 58333
$ws.Range("K141").Value = 20000
$ws.Range("L141").Value = 58333
$ws.Range("M141").Value = -14820
$ws.Range("N141").Value = -68693

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 315.77777
$ws.Range("I23").Value = 5.5
$ws.Range("J23").Value = 354.5625
$ws.Range("K23").Value = 16.5
$ws.Range("L23").Value = 1063.6875
$ws.Range("M23").Value = 218.5
$ws.Range("N23").Value = -1533.6875

$ws.Range("H26").Value = 65.5
$ws.Range("I26").Value = 65.5
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 196.5
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 91.5
$ws.Range("N26").ClearContents()

$ws.Range("H122").Value = 710.26666
$ws.Range("I122").Value = 500.1111
$ws.Range("J122").Value = 1025.5
$ws.Range("K122").Value = 4500.9999
$ws.Range("L122").Value = 9229.5
$ws.Range("M122").Value = -2050.9999
$ws.Range("N122").Value = -14129.5

$ws.Range("H132").Value = 2474.3125
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2474.3125
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 22268.8125
$ws.Range("N132").Value = -27328.8125
$ws.Range("M132").ClearContents()

$ws.Range("H134").Value = 4421.387
$ws.Range("I134").Value = 2639.1667
$ws.Range("J134").Value = 5547
$ws.Range("K134").Value = 7917.500100000001
$ws.Range("L134").Value = 16641
$ws.Range("M134").Value = -2847.500100000001
$ws.Range("N134").Value = -26781

$ws.Range("H138").Value = 2787.1052
$ws.Range("I138").Value = 906.5
$ws.Range("J138").Value = 3655.077
$ws.Range("K138").Value = 2719.5
$ws.Range("L138").Value = 10965.231
$ws.Range("M138").Value = 2420.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 9274.166999999999
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 9274.166999999999
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 9274.166999999999
$ws.Range("N109").Value = -11354.167

$ws.Range("H123").Value = 10138.733
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 10138.733
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 10138.733
$ws.Range("N123").Value = -15038.733

$ws.Range("H132").Value = 1784.0416
$ws.Range("I132").Value = 1296.28
$ws.Range("J132").Value = 2314.2173
$ws.Range("K132").Value = 3888.84
$ws.Range("L132").Value = 6942.651899999999
$ws.Range("M132").Value = -1358.84
$ws.Range("N132").Value = -12002.6519

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1145.2174
$ws.Range("I46").Value = 921.4286
$ws.Range("J46").Value = 1493.3334
$ws.Range("K46").Value = 921.4286
$ws.Range("L46").Value = 1493.3334
$ws.Range("M46").Value = -733.4286
$ws.Range("N46").Value = -1869.3334

$ws.Range("H55").Value = 799.1429000000001
$ws.Range("I55").Value = 395
$ws.Range("J55").Value = 960.8
$ws.Range("K55").Value = 395
$ws.Range("L55").Value = 960.8
$ws.Range("M55").Value = -222
$ws.Range("N55").Value = -1306.8

$ws.Range("H138").Value = 40000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 40000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 40000
$ws.Range("N138").Value = -50280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4260
$ws.Range("I62").Value = 3833.3333
$ws.Range("J62").Value = 4900
$ws.Range("K62").Value = 3833.3333
$ws.Range("L62").Value = 4900
$ws.Range("M62").Value = -3209.3333
$ws.Range("N62").Value = -6148

$ws.Range("H65").Value = 4260
$ws.Range("I65").Value = 3833.3333
$ws.Range("J65").Value = 4900
$ws.Range("K65").Value = 19166.6665
$ws.Range("L65").Value = 24500
$ws.Range("M65").Value = -16046.6665
$ws.Range("N65").Value = -30740

$ws.Range("H136").Value = 1779.9259
$ws.Range("I136").Value = 1877.4166
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 5632.2498
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -3082.2498
